$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column (H1), matching the formatting of the
# existing header cells (e.g. G1 "sum")
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Add the corresponding data value for row 2
$ws.Range("H2").Value = 0
